# Update the 25 two-digit multiplication "problem = answer" cells in the
# table to regenerate the worksheet content, matching commit
# "Update master to output generated at aa3dc9e".
$d = $word.ActiveDocument
$d.Content.Find.Execute("59×74=4366", $false, $false, $false, $false, $false, $true, 1, $false, "34×67=2278", 2) | Out-Null
$d.Content.Find.Execute("70×43=3010", $false, $false, $false, $false, $false, $true, 1, $false, "56×78=4368", 2) | Out-Null
$d.Content.Find.Execute("85×29=2465", $false, $false, $false, $false, $false, $true, 1, $false, "66×52=3432", 2) | Out-Null
$d.Content.Find.Execute("21×37=777", $false, $false, $false, $false, $false, $true, 1, $false, "65×97=6305", 2) | Out-Null
$d.Content.Find.Execute("43×70=3010", $false, $false, $false, $false, $false, $true, 1, $false, "61×93=5673", 2) | Out-Null
$d.Content.Find.Execute("88×53=4664", $false, $false, $false, $false, $false, $true, 1, $false, "99×17=1683", 2) | Out-Null
$d.Content.Find.Execute("37×86=3182", $false, $false, $false, $false, $false, $true, 1, $false, "69×38=2622", 2) | Out-Null
$d.Content.Find.Execute("96×78=7488", $false, $false, $false, $false, $false, $true, 1, $false, "51×85=4335", 2) | Out-Null
$d.Content.Find.Execute("52×86=4472", $false, $false, $false, $false, $false, $true, 1, $false, "40×83=3320", 2) | Out-Null
$d.Content.Find.Execute("82×43=3526", $false, $false, $false, $false, $false, $true, 1, $false, "62×13=806", 2) | Out-Null
$d.Content.Find.Execute("99×91=9009", $false, $false, $false, $false, $false, $true, 1, $false, "29×18=522", 2) | Out-Null
$d.Content.Find.Execute("29×45=1305", $false, $false, $false, $false, $false, $true, 1, $false, "24×96=2304", 2) | Out-Null
$d.Content.Find.Execute("60×28=1680", $false, $false, $false, $false, $false, $true, 1, $false, "35×57=1995", 2) | Out-Null
$d.Content.Find.Execute("40×39=1560", $false, $false, $false, $false, $false, $true, 1, $false, "85×48=4080", 2) | Out-Null
$d.Content.Find.Execute("63×89=5607", $false, $false, $false, $false, $false, $true, 1, $false, "60×51=3060", 2) | Out-Null
$d.Content.Find.Execute("81×88=7128", $false, $false, $false, $false, $false, $true, 1, $false, "34×66=2244", 2) | Out-Null
$d.Content.Find.Execute("64×97=6208", $false, $false, $false, $false, $false, $true, 1, $false, "50×63=3150", 2) | Out-Null
$d.Content.Find.Execute("64×98=6272", $false, $false, $false, $false, $false, $true, 1, $false, "41×96=3936", 2) | Out-Null
$d.Content.Find.Execute("39×71=2769", $false, $false, $false, $false, $false, $true, 1, $false, "64×67=4288", 2) | Out-Null
$d.Content.Find.Execute("17×59=1003", $false, $false, $false, $false, $false, $true, 1, $false, "61×35=2135", 2) | Out-Null
$d.Content.Find.Execute("37×91=3367", $false, $false, $false, $false, $false, $true, 1, $false, "63×54=3402", 2) | Out-Null
$d.Content.Find.Execute("26×86=2236", $false, $false, $false, $false, $false, $true, 1, $false, "24×72=1728", 2) | Out-Null
$d.Content.Find.Execute("79×67=5293", $false, $false, $false, $false, $false, $true, 1, $false, "13×24=312", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $false, $false, $false, $false, $false, $true, 1, $false, "62×44=2728", 2) | Out-Null
$d.Content.Find.Execute("13×17=221", $false, $false, $false, $false, $false, $true, 1, $false, "62×80=4960", 2) | Out-Null
